$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "optimised." -> "supervised."
#
# The paragraph also carries Word's auto-managed "_GoBack" bookmark (marking
# the last edit location). In the target revision that bookmark sits right
# between "supervised" and the trailing "." instead of in the following
# empty paragraph, so after changing the text we relocate the bookmark to
# match.
# ---------------------------------------------------------------------------
$rOptim = $d.Content.Duplicate
$rOptim.Find.MatchCase = $true
$rOptim.Find.Execute("optimised.")
if ($rOptim.Find.Found) {
    $optimStart = $rOptim.Start
    $rOptim.Text = "supervised."

    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }
    $bmPos = $optimStart + 10
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ---------------------------------------------------------------------------
# Change 2: expand the conclusion sentence with an additional clause.
#
# ", ... outputs and take decisions on those outcomes. " becomes
# ", ... outputs, take decisions on those outcomes and finally limit power
#   consumption in a given month. "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "outputs and take decisions on those outcomes.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "outputs, take decisions on those outcomes and finally limit power consumption in a given month.", `
    2)
